$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
Write-Host $ws.Name
